$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 1.8
$ws.Range("G2").Value = 1.82
$ws.Range("H2").Value = 5.2
$ws.Range("I2").Value = 5.3
$ws.Range("J2").Value = 3.85
$ws.Range("K2").Value = 3.9
$ws.Range("L2").Value = 1.46
$ws.Range("N2").Value = 3.6
$ws.Range("O2").Value = 1.37
$ws.Range("P2").Value = 1.86
$ws.Range("Q2").Value = 2.12
$ws.Range("R2").Value = 1.33
$ws.Range("S2").Value = 3.9
$ws.Range("T2").Value = 2
$ws.Range("U2").Value = 1.95
$ws.Range("V2").Value = 1.23
$ws.Range("W2").Value = 2.2
$ws.Range("Y2").Value = 16
$ws.Range("Z2").Value = 38
$ws.Range("AA2").Value = 150
$ws.Range("AB2").Value = 7.8
$ws.Range("AC2").Value = 8.2
$ws.Range("AD2").Value = 19.5
$ws.Range("AE2").Value = 75
$ws.Range("AF2").Value = 9.6
$ws.Range("AG2").Value = 10.5
$ws.Range("AH2").Value = 22
$ws.Range("AI2").Value = 85
$ws.Range("AJ2").Value = 18.5
$ws.Range("AK2").Value = 19
$ws.Range("AN2").Value = 13.5
$ws.Range("AO2").Value = 100

$ws.Range("F3").Value = 1.39
$ws.Range("G3").Value = 1.47
$ws.Range("I3").Value = 9.6
$ws.Range("J3").Value = 4.7
$ws.Range("K3").Value = 7.2
$ws.Range("L3").Value = 1.25
$ws.Range("M3").Value = 1.04
$ws.Range("N3").Value = 5.5
$ws.Range("P3").Value = 2.74
$ws.Range("Q3").Value = 1.48
$ws.Range("R3").Value = 1.71
$ws.Range("S3").Value = 2.22
$ws.Range("T3").Value = 1.7
$ws.Range("U3").Value = 2.14
$ws.Range("V3").Value = 1.12
$ws.Range("W3").Value = 3.1
$ws.Range("AB3").Value = 990
$ws.Range("AG3").Value = 990
$ws.Range("AJ3").Value = 900

$ws.Range("F4").Value = 1.95
$ws.Range("G4").Value = 2.04
$ws.Range("H4").Value = 4.4
$ws.Range("I4").Value = 5.1
$ws.Range("J4").Value = 3.3
$ws.Range("M4").Value = 1.11
$ws.Range("N4").Value = 2.92
$ws.Range("O4").Value = 1.47
$ws.Range("P4").Value = 1.63
$ws.Range("Q4").Value = 2.38
$ws.Range("R4").Value = 1.23
$ws.Range("S4").Value = 4.7
$ws.Range("T4").Value = 2.02
$ws.Range("U4").Value = 1.83
$ws.Range("V4").Value = 1.24
$ws.Range("W4").Value = 1.96
$ws.Range("X4").Value = 10.5
$ws.Range("Y4").Value = 13.5
$ws.Range("Z4").Value = 1000
$ws.Range("AH4").Value = 42
$ws.Range("AJ4").Value = 32
$ws.Range("AK4").Value = 38
$ws.Range("AL4").Value = 130

$ws.Range("F5").Value = 2.36
$ws.Range("G5").Value = 2.52
$ws.Range("H5").Value = 3.2
$ws.Range("I5").Value = 3.5
$ws.Range("J5").Value = 3.4
$ws.Range("K5").Value = 3.55
$ws.Range("M5").Value = 1.1
$ws.Range("N5").Value = 3.25
$ws.Range("O5").Value = 1.42
$ws.Range("P5").Value = 1.74
$ws.Range("Q5").Value = 2.22
$ws.Range("S5").Value = 4.1
$ws.Range("T5").Value = 1.83
$ws.Range("V5").Value = 1.41
$ws.Range("W5").Value = 1.62
$ws.Range("Y5").Value = 14.5
$ws.Range("AA5").Value = 1000
$ws.Range("AB5").Value = 22
$ws.Range("AE5").Value = 50
$ws.Range("AF5").Value = 36
$ws.Range("AH5").Value = 60

$ws.Range("F6").Value = 1.69
$ws.Range("G6").Value = 1.81
$ws.Range("H6").Value = 4.7
$ws.Range("I6").Value = 5.5
$ws.Range("J6").Value = 4.1
$ws.Range("K6").Value = 4.8
$ws.Range("M6").Value = 1.05
$ws.Range("N6").Value = 4.7
$ws.Range("P6").Value = 2.28
$ws.Range("Q6").Value = 1.67
$ws.Range("R6").Value = 1.51
$ws.Range("S6").Value = 2.68
$ws.Range("T6").Value = 1.69
$ws.Range("U6").Value = 2.16
$ws.Range("W6").Value = 2.22
$ws.Range("AC6").Value = 970

$ws.Range("H7").Value = 1.81
$ws.Range("I7").Value = 1.83
$ws.Range("J7").Value = 3.6
$ws.Range("K7").Value = 3.65
$ws.Range("L7").Value = 1.51
$ws.Range("N7").Value = 3.25
$ws.Range("O7").Value = 1.42
$ws.Range("P7").Value = 1.75
$ws.Range("Q7").Value = 2.3
$ws.Range("R7").Value = 1.27
$ws.Range("S7").Value = 4.3
$ws.Range("T7").Value = 2.1
$ws.Range("U7").Value = 1.86
$ws.Range("V7").Value = 2.2
$ws.Range("W7").Value = 1.2
$ws.Range("Y7").Value = 7.2
$ws.Range("Z7").Value = 9.4
$ws.Range("AB7").Value = 16.5
$ws.Range("AI7").Value = 46
$ws.Range("AJ7").Value = 150
$ws.Range("AL7").Value = 100
$ws.Range("AM7").Value = 160
$ws.Range("AN7").Value = 140
$ws.Range("AO7").Value = 15

$ws.Range("F8").Value = 2.8
$ws.Range("H8").Value = 2.4
$ws.Range("I8").Value = 2.56
$ws.Range("J8").Value = 3.7
$ws.Range("K8").Value = 3.95
$ws.Range("M8").Value = 1.06
$ws.Range("P8").Value = 2.12
$ws.Range("Q8").Value = 1.78
$ws.Range("S8").Value = 2.94
$ws.Range("T8").Value = 1.64
$ws.Range("U8").Value = 2.28
$ws.Range("V8").Value = 1.64
$ws.Range("X8").Value = 18.5
$ws.Range("Y8").Value = 13
$ws.Range("Z8").Value = 18.5
$ws.Range("AA8").Value = 36
$ws.Range("AB8").Value = 14.5
$ws.Range("AC8").Value = 8.8
$ws.Range("AE8").Value = 26

$ws.Range("F9").Value = 5.7
$ws.Range("G9").Value = 5.8
$ws.Range("H9").Value = 1.82
$ws.Range("I9").Value = 1.84
$ws.Range("J9").Value = 3.55
$ws.Range("K9").Value = 3.6
$ws.Range("L9").Value = 1.55
$ws.Range("N9").Value = 3
$ws.Range("O9").Value = 1.48
$ws.Range("P9").Value = 1.67
$ws.Range("Q9").Value = 2.46
$ws.Range("S9").Value = 4.9
$ws.Range("V9").Value = 2.18
$ws.Range("W9").Value = 1.21
$ws.Range("X9").Value = 9.8
$ws.Range("Z9").Value = 9.2
$ws.Range("AB9").Value = 15
$ws.Range("AC9").Value = 8
$ws.Range("AH9").Value = 25
$ws.Range("AJ9").Value = 160
$ws.Range("AL9").Value = 110
$ws.Range("AO9").Value = 17

$ws.Range("L10").Value = 1.43
$ws.Range("P10").Value = 1.94
$ws.Range("T10").Value = 1.95
$ws.Range("U10").Value = 2.02
$ws.Range("X10").Value = 14
$ws.Range("AC10").Value = 8.4

$ws.Range("F11").Value = 3.15
$ws.Range("G11").Value = 3.25
$ws.Range("N11").Value = 3.7
$ws.Range("P11").Value = 1.9
$ws.Range("Q11").Value = 2.06
$ws.Range("W11").Value = 1.44
$ws.Range("AC11").Value = 7.8
$ws.Range("AI11").Value = 42
$ws.Range("AM11").Value = 330
$ws.Range("AN11").Value = 34

$ws.Range("G12").Value = 1.16
$ws.Range("J12").Value = 9.8
$ws.Range("K12").Value = 12
$ws.Range("N12").Value = 8.4
$ws.Range("P12").Value = 3.5
$ws.Range("Q12").Value = 1.32
$ws.Range("R12").Value = 2.04
$ws.Range("S12").Value = 1.8
$ws.Range("U12").Value = 1.7
$ws.Range("AC12").Value = 26
$ws.Range("AH12").Value = 330

$ws.Range("J13").Value = 3.5
$ws.Range("K13").Value = 3.55
$ws.Range("L13").Value = 1.47
$ws.Range("M13").Value = 1.08
$ws.Range("N13").Value = 3.65
$ws.Range("O13").Value = 1.36
$ws.Range("P13").Value = 1.9
$ws.Range("Q13").Value = 2.1
$ws.Range("R13").Value = 1.34
$ws.Range("S13").Value = 3.85
$ws.Range("T13").Value = 1.84
$ws.Range("U13").Value = 2.12
$ws.Range("X13").Value = 13.5
$ws.Range("Y13").Value = 9.6
$ws.Range("Z13").Value = 14
$ws.Range("AB13").Value = 12
$ws.Range("AC13").Value = 7.4
$ws.Range("AE13").Value = 25
$ws.Range("AF13").Value = 21
$ws.Range("AI13").Value = 38
$ws.Range("AJ13").Value = 55
$ws.Range("AL13").Value = 48
$ws.Range("AM13").Value = 85
$ws.Range("AN13").Value = 40
$ws.Range("AO13").Value = 24

$ws.Range("G14").Value = 1.57
$ws.Range("H14").Value = 7.4
$ws.Range("J14").Value = 4.4
$ws.Range("O14").Value = 1.22
$ws.Range("S14").Value = 2.78
$ws.Range("T14").Value = 1.83
$ws.Range("W14").Value = 2.74
$ws.Range("X14").Value = 21
$ws.Range("Y14").Value = 29
$ws.Range("Z14").Value = 150
$ws.Range("AC14").Value = 10
$ws.Range("AD14").Value = 26
$ws.Range("AE14").Value = 100
$ws.Range("AG14").Value = 10
$ws.Range("AH14").Value = 19.5
$ws.Range("AI14").Value = 90
$ws.Range("AO14").Value = 110

